$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime the style used by the highlighted "use case" rows (fillId 4 / style 7),
# by copying an existing cell's formatting onto the new C47:C51 cells.
$ws.Range("C29").Copy()
$ws.Range("C47:C51").PasteSpecial(-4122)  # xlPasteFormats

# --- New section header: "Sviluppo 4 casi d'uso " (row 47) ---
$ws.Range("C47").Value = "Sviluppo 4 casi d'uso "

# --- Fill in the C column for the 4 use cases, in the order the author typed them ---
$ws.Range("C50").Value = "Autenticazione utente"
$ws.Range("C51").Value = "Assegna ruolo utente"

# --- Status column ---
$ws.Range("D48").Value = "implementato"
$ws.Range("D49").Value = "implementato"

$ws.Range("D29").Value = "completato"
$ws.Range("D30").Value = "completato"
$ws.Range("D31").Value = "completato"
$ws.Range("D29:D31").Font.Bold = $true

$ws.Range("C48").Value = "Consultare dettaglio POI"
$ws.Range("C49").Value = "Registrazione nuovo utente"

$ws.Range("F48").Value = "dettagliato con diagrammi"
$ws.Range("F49").Value = "dettagliato con diagrammi"
$ws.Range("F50").Value = "dettagliato con diagrammi"
$ws.Range("F51").Value = "dettagliato con diagrammi"

$ws.Range("D50").Value = "da implementare"
$ws.Range("D51").Value = "da implementare"

# --- Update selection to reflect scrolling to the new rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$ws.Range("D37").Select()
